$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update version string (Property "Version" row, B3)
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update date (Property "Date" row, B8)
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row for "Jurisdiction" right after "Contact" (row 10) and
# before "Description" (old row 11), pushing the remaining property rows
# down by one.
$ws.Rows.Item(11).Insert()

# Copy formatting from the row above so the new row matches the rest of
# the table (borders / alignment / style).
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
